# Update the split-purchase dataset: Item / Item Category / Vendor columns
# for the 8 data rows (rows 2-9) on the "PR Document" sheet, add the
# AutoFilter defined name for the used range, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Item (column D) -------------------------------------------------
# Set column-by-column so new shared strings are interned in the same
# order as the source data (TV, Stereo System, Audio Visual System,
# Vendor_1..4).
$ws.Range("D2").Value = "TV"
$ws.Range("D3").Value = "TV"
$ws.Range("D4").Value = "Stereo System"
$ws.Range("D5").Value = "Stereo System"
$ws.Range("D6").Value = "Stereo System"
$ws.Range("D7").Value = "TV"
$ws.Range("D8").Value = "Stereo System"
$ws.Range("D9").Value = "Stereo System"

# --- Item Category (column E) ----------------------------------------
$ws.Range("E2").Value = "Audio Visual System"
$ws.Range("E3").Value = "Audio Visual System"
$ws.Range("E4").Value = "Audio Visual System"
$ws.Range("E5").Value = "Audio Visual System"
$ws.Range("E6").Value = "Audio Visual System"
$ws.Range("E7").Value = "Audio Visual System"
$ws.Range("E8").Value = "Audio Visual System"
$ws.Range("E9").Value = "Audio Visual System"

# --- Vendor (column F) -------------------------------------------------
$ws.Range("F2").Value = "Vendor_1"
$ws.Range("F3").Value = "Vendor_1"
$ws.Range("F4").Value = "Vendor_2"
$ws.Range("F5").Value = "Vendor_3"
$ws.Range("F6").Value = "Vendor_4"
$ws.Range("F7").Value = "Vendor_1"
$ws.Range("F8").Value = "Vendor_2"
$ws.Range("F9").Value = "Vendor_2"

# --- AutoFilter defined name (sheet-scoped, hidden) --------------------
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='PR Document'!`$A`$1:`$G`$9")
$filterName.Visible = $false

# --- Selection ----------------------------------------------------------
$ws.Range("F19").Select() | Out-Null
